$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.564
$ws.Range("A6").Value = -22.255
$ws.Range("A7").Value = -19.749
$ws.Range("C7").Value = -12.489
$ws.Range("C12").Value = -10.782
$ws.Range("D13").Value = -8.029999999999998
$ws.Range("D14").Value = -7.787000000000001
$ws.Range("C15").Value = -13.286
$ws.Range("A16").Value = -21.774
$ws.Range("D16").Value = -8.597999999999999
$ws.Range("D19").Value = -8.266
$ws.Range("A20").Value = -19.765
$ws.Range("C20").Value = -12.366
$ws.Range("C21").Value = -12.293
$ws.Range("C22").Value = -12.8
$ws.Range("D22").Value = -7.783000000000001
$ws.Range("C23").Value = -12.223
$ws.Range("A28").Value = -21.942
$ws.Range("A29").Value = -21.426
$ws.Range("C29").Value = -12.305
$ws.Range("A32").Value = -21.767
$ws.Range("C34").Value = -11.955
$ws.Range("D36").Value = -8.189
$ws.Range("A40").Value = -20.097
$ws.Range("C42").Value = -12.574
$ws.Range("C43").Value = -13.159
$ws.Range("C44").Value = -13.216
$ws.Range("C45").Value = -13.054
$ws.Range("A46").Value = -21.877
$ws.Range("C46").Value = -13.572
$ws.Range("D46").Value = -8.682
$ws.Range("C50").Value = -13.625
$ws.Range("D50").Value = -8.097
$ws.Range("A51").Value = -21.684
$ws.Range("C51").Value = -11.06
$ws.Range("A52").Value = -21.907
$ws.Range("A57").Value = -22.247
$ws.Range("A59").Value = -22.29
$ws.Range("A62").Value = -22.158
$ws.Range("A66").Value = -21.684
$ws.Range("C66").Value = -11.405
$ws.Range("C67").Value = -11.3
$ws.Range("A73").Value = -20.597
$ws.Range("A74").Value = -21.244
$ws.Range("C79").Value = -11.889
$ws.Range("C84").Value = -14.098
$ws.Range("A92").Value = -21.626
$ws.Range("C92").Value = -11.472
$ws.Range("D95").Value = -7.854000000000001
$ws.Range("C97").Value = -12.592
$ws.Range("D97").Value = -8.561
$ws.Range("A100").Value = -22.063
